# "fixed obs data for graphs"
#
# 1) ObservedSoil sheet: the SimulationName column (A) used the wrong
#    simulation names (Gm / Gm48N / MakokaN0 / MakokaN48). Rename them to
#    the correct Makoka*Maize* names so they line up with the simulations
#    used for the graphs.
# 2) ObservedGliricidia sheet: a handful of sample dates were off by one
#    day - correct them.

$wb = $excel.ActiveWorkbook

$wsSoil = $wb.Worksheets.Item("ObservedSoil")
$wsGlir = $wb.Worksheets.Item("ObservedGliricidia")

# --- ObservedSoil: fix the SimulationName values in column A -------------
$wsSoil.Range("A2:A17").Value  = "MakokaN0MaizeSole"
$wsSoil.Range("A18:A33").Value = "MakokaN48MaizeSole"
$wsSoil.Range("A34:A49").Value = "MakokaN0MaizeGliricidia"
$wsSoil.Range("A50:A65").Value = "MakokaN48MaizeGliricidia"

# Move selection back to the top of the sheet once the edits are done.
$wsSoil.Range("B1").Select() | Out-Null

# --- ObservedGliricidia: correct a few sample dates -----------------------
$wsGlir.Activate() | Out-Null
$wsGlir.Range("C37").Value = 35289
$wsGlir.Range("C43").Value = 35245
$wsGlir.Range("C46").Value = 35291
$wsGlir.Range("C47").Value = 35347
$wsGlir.Range("C49").Value = 35399

$wsGlir.Range("C38").Select() | Out-Null
